# Auto-generated edit script: apply cell value updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("AR3").Value = 4.8
$ws.Range("AS3").Value = 1.19
# Row 4
$ws.Range("AR4").Value = 5
$ws.Range("AS4").Value = 1.18
# Row 5
$ws.Range("G5").Value = 2.63
$ws.Range("I5").Value = 3.1
$ws.Range("AJ5").Value = 7
$ws.Range("AN5").Value = 34
$ws.Range("AR5").Value = 4.4
$ws.Range("AS5").Value = 1.22
# Row 6
$ws.Range("AR6").Value = 4.1
$ws.Range("AS6").Value = 1.24
# Row 7
$ws.Range("AR7").Value = 5.2
$ws.Range("AS7").Value = 1.17
# Row 8
$ws.Range("G8").Value = 2.25
$ws.Range("I8").Value = 2.7
$ws.Range("J8").Value = 2.88
$ws.Range("AK8").Value = 17
$ws.Range("AL8").Value = 11
# Row 10
$ws.Range("G10").Value = 2.8
$ws.Range("I10").Value = 2.3
$ws.Range("J10").Value = 3.4
$ws.Range("AJ10").Value = 9.5
$ws.Range("AN10").Value = 17
# Row 11
$ws.Range("G11").Value = 1.62
$ws.Range("H11").Value = 3.8
$ws.Range("I11").Value = 5.5
$ws.Range("M11").Value = 1.07
$ws.Range("N11").Value = 9
$ws.Range("W11").Value = 2.2
$ws.Range("X11").Value = 1.62
$ws.Range("AB11").Value = 12
$ws.Range("AG11").Value = 21
$ws.Range("AJ11").Value = 11
$ws.Range("AK11").Value = 26
$ws.Range("AP11").Value = 1.66
$ws.Range("AQ11").Value = 2.19
$ws.Range("AR11").Value = 3.3
$ws.Range("AS11").Value = 1.32
# Row 12
$ws.Range("G12").Value = 2.3
$ws.Range("H12").Value = 2.8
$ws.Range("J12").Value = 2.95
$ws.Range("K12").Value = 1.9
$ws.Range("L12").Value = 4.05
$ws.Range("M12").Value = 1.12
$ws.Range("N12").Value = 5.3
$ws.Range("O12").Value = 1.5
$ws.Range("P12").Value = 2.42
$ws.Range("Q12").Value = 2.5
$ws.Range("S12").Value = 4.4
$ws.Range("U12").Value = 1.53
$ws.Range("V12").Value = 2.35
$ws.Range("W12").Value = 2.02
$ws.Range("X12").Value = 1.72
$ws.Range("Z12").Value = 9.75
$ws.Range("AE12").Value = 5.3
$ws.Range("AJ12").Value = 8
$ws.Range("AK12").Value = 17.5
$ws.Range("AL12").Value = 12
$ws.Range("AN12").Value = 37
$ws.Range("AO12").Value = 50
# Row 15
$ws.Range("G15").Value = 1.4
$ws.Range("I15").Value = 6.5
$ws.Range("AA15").Value = 8.5
$ws.Range("AB15").Value = 9.5
$ws.Range("AD15").Value = 26
$ws.Range("AE15").Value = 13
$ws.Range("AF15").Value = 9.5
$ws.Range("AJ15").Value = 17
$ws.Range("AK15").Value = 34
$ws.Range("AL15").Value = 19
# Row 18
$ws.Range("G18").Value = 3.4
$ws.Range("H18").Value = 3.3
$ws.Range("I18").Value = 2.1
$ws.Range("L18").Value = 2.88
$ws.Range("Q18").Value = 2.15
$ws.Range("R18").Value = 1.67
$ws.Range("Y18").Value = 9
$ws.Range("Z18").Value = 17
$ws.Range("AB18").Value = 41
$ws.Range("AF18").Value = 6.5
$ws.Range("AK18").Value = 9.5
$ws.Range("AL18").Value = 9
$ws.Range("AM18").Value = 19
$ws.Range("AN18").Value = 19
# Row 19
$ws.Range("G19").Value = 2.38
$ws.Range("H19").Value = 3.3
$ws.Range("I19").Value = 2.9
$ws.Range("J19").Value = 3.1
$ws.Range("Y19").Value = 7.5
$ws.Range("Z19").Value = 11
$ws.Range("AA19").Value = 9.5
$ws.Range("AK19").Value = 15
# Row 20
$ws.Range("G20").Value = 1.53
$ws.Range("H20").Value = 4.33
$ws.Range("I20").Value = 5.5
$ws.Range("L20").Value = 5.5
$ws.Range("AE20").Value = 15
$ws.Range("AJ20").Value = 15
$ws.Range("AM20").Value = 51
# Row 21
$ws.Range("G21").Value = 2.55
$ws.Range("H21").Value = 3.25
$ws.Range("I21").Value = 2.7
$ws.Range("J21").Value = 3.25
$ws.Range("M21").Value = 1.06
$ws.Range("N21").Value = 10
$ws.Range("O21").Value = 1.3
$ws.Range("P21").Value = 3.4
$ws.Range("Q21").Value = 2.03
$ws.Range("R21").Value = 1.83
$ws.Range("S21").Value = 3.5
$ws.Range("T21").Value = 1.29
$ws.Range("W21").Value = 1.75
$ws.Range("X21").Value = 2
$ws.Range("Y21").Value = 8.5
$ws.Range("Z21").Value = 13
$ws.Range("AA21").Value = 10
$ws.Range("AB21").Value = 26
$ws.Range("AC21").Value = 21
$ws.Range("AE21").Value = 10
$ws.Range("AF21").Value = 6
$ws.Range("AH21").Value = 51
$ws.Range("AJ21").Value = 9
$ws.Range("AK21").Value = 13
$ws.Range("AL21").Value = 10
$ws.Range("AM21").Value = 26
# Row 22
$ws.Range("H22").Value = 3.8
$ws.Range("J22").Value = 2.3
$ws.Range("K22").Value = 2.2
$ws.Range("M22").Value = 1.06
$ws.Range("N22").Value = 10
$ws.Range("O22").Value = 1.33
$ws.Range("P22").Value = 3.25
$ws.Range("Q22").Value = 2.05
$ws.Range("R22").Value = 1.75
$ws.Range("S22").Value = 3.5
$ws.Range("T22").Value = 1.29
$ws.Range("W22").Value = 1.95
$ws.Range("X22").Value = 1.8
$ws.Range("Y22").Value = 6.5
$ws.Range("AD22").Value = 29
$ws.Range("AE22").Value = 9.5
$ws.Range("AI22").Value = 401
$ws.Range("AJ22").Value = 12
# Row 23
$ws.Range("J23").Value = 2
$ws.Range("W23").Value = 2.2
$ws.Range("X23").Value = 1.62
$ws.Range("Z23").Value = 6
# Row 24
$ws.Range("G24").Value = 1.95
$ws.Range("H24").Value = 3.1
$ws.Range("I24").Value = 4.33
$ws.Range("J24").Value = 2.75
$ws.Range("K24").Value = 1.91
$ws.Range("M24").Value = 1.11
$ws.Range("N24").Value = 6.5
$ws.Range("O24").Value = 1.53
$ws.Range("P24").Value = 2.38
$ws.Range("Q24").Value = 2.7
$ws.Range("R24").Value = 1.44
$ws.Range("S24").Value = 5.5
$ws.Range("T24").Value = 1.14
$ws.Range("U24").Value = 1.62
$ws.Range("V24").Value = 2.2
$ws.Range("W24").Value = 2.25
$ws.Range("X24").Value = 1.57
$ws.Range("Y24").Value = 5
$ws.Range("AA24").Value = 9.5
$ws.Range("AB24").Value = 17
$ws.Range("AC24").Value = 21
$ws.Range("AE24").Value = 6
$ws.Range("AJ24").Value = 8.5
$ws.Range("AL24").Value = 17
$ws.Range("AP24").Value = 2.03
$ws.Range("AQ24").Value = 1.83
# Row 25
$ws.Range("Q25").Value = 1.95
$ws.Range("R25").Value = 1.9
$ws.Range("S25").Value = 3.25
$ws.Range("T25").Value = 1.33
# Row 26
$ws.Range("N26").Value = 7.5
$ws.Range("Q26").Value = 2.35
$ws.Range("R26").Value = 1.57
$ws.Range("S26").Value = 4.33
$ws.Range("T26").Value = 1.2
$ws.Range("U26").Value = 1.53
$ws.Range("V26").Value = 2.38
$ws.Range("W26").Value = 2
$ws.Range("X26").Value = 1.75
$ws.Range("Y26").Value = 8
$ws.Range("AE26").Value = 7.5
$ws.Range("AG26").Value = 17
$ws.Range("AI26").Value = 401
$ws.Range("AJ26").Value = 6.5
# Row 28
$ws.Range("G28").Value = 2.37
$ws.Range("I28").Value = 2.75
$ws.Range("J28").Value = 2.87
$ws.Range("K28").Value = 2.15
$ws.Range("L28").Value = 3.25
$ws.Range("Q28").Value = 1.52
$ws.Range("R28").Value = 2.2
$ws.Range("S28").Value = 2.2
$ws.Range("T28").Value = 1.53
$ws.Range("W28").Value = 1.42
$ws.Range("X28").Value = 2.47
$ws.Range("Y28").Value = 12
$ws.Range("Z28").Value = 16
$ws.Range("AA28").Value = 9.25
$ws.Range("AB28").Value = 29
$ws.Range("AC28").Value = 17
$ws.Range("AD28").Value = 19
$ws.Range("AE28").Value = 14
$ws.Range("AF28").Value = 6.8
$ws.Range("AG28").Value = 10.25
$ws.Range("AH28").Value = 32
$ws.Range("AI28").Value = 175
$ws.Range("AJ28").Value = 12.5
$ws.Range("AK28").Value = 17.5
$ws.Range("AL28").Value = 10
$ws.Range("AN28").Value = 21
$ws.Range("AO28").Value = 21
# Row 29
$ws.Range("G29").Value = 1.45
$ws.Range("J29").Value = 1.95
$ws.Range("K29").Value = 2.4
$ws.Range("M29").Value = 1.03
$ws.Range("N29").Value = 15
$ws.Range("Q29").Value = 1.7
$ws.Range("R29").Value = 2.1
# Row 30
$ws.Range("N30").Value = 9
$ws.Range("O30").Value = 1.36
$ws.Range("P30").Value = 3
$ws.Range("Q30").Value = 2.15
$ws.Range("R30").Value = 1.67
# Row 33
$ws.Range("G33").Value = 2.75
$ws.Range("M33").Value = 1.04
$ws.Range("N33").Value = 9
$ws.Range("R33").Value = 1.75
$ws.Range("W33").Value = 1.8
$ws.Range("X33").Value = 1.91
$ws.Range("AC33").Value = 23
$ws.Range("AE33").Value = 9
$ws.Range("AG33").Value = 15
$ws.Range("AH33").Value = 51
$ws.Range("AK33").Value = 12
